$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.954.88'
$ws.Range("E2").Value = '  -0.16%  '
$ws.Range("D3").Value = '1.996.46'
$ws.Range("E3").Value = '  -0.85%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '255.28'
$ws.Range("E5").Value = '  +3.91%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.615'
$ws.Range("E6").Value = '  -1.09%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '55.27'
$ws.Range("E8").Value = '  -6.39%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.375'
$ws.Range("E9").Value = '  -3.90%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0763'
$ws.Range("E10").Value = '  -4.68%  '
$ws.Range("E11").Value = '  -2.61%  '
$ws.Range("D12").Value = '2.292.63'
$ws.Range("E12").Value = '  -0.69%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.11'
$ws.Range("E13").Value = '  -5.62%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.17'
$ws.Range("E14").Value = '  -4.09%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.787'
$ws.Range("E15").Value = '  -6.48%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.15'
$ws.Range("E16").Value = '  -4.65%  '
$ws.Range("D17").Value = '1.995.69'
$ws.Range("E17").Value = '  -0.85%  '
$ws.Range("D18").Value = '36.831.01'
$ws.Range("E18").Value = '  -0.31%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '70.35'
$ws.Range("E19").Value = '  +0.50%  '
$ws.Range("D20").Value = '0.0₃0823'
$ws.Range("E20").Value = '  -4.00%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '233.33'
$ws.Range("E21").Value = '  +1.78%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.04'
$ws.Range("E22").Value = '  -2.85%  '
$ws.Range("E23").Value = '  +0.19%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.51'
$ws.Range("E24").Value = '  -1.89%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.36'
$ws.Range("E25").Value = '  +0.19%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '163.98'
$ws.Range("E26").Value = '  +0.34%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.82'
$ws.Range("E27").Value = '  -5.30%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.38'
$ws.Range("E28").Value = '  -1.54%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.124'
$ws.Range("E29").Value = '  -9.02%  '
$ws.Range("E30").Value = '  -3.72%  '
$ws.Range("E31").Value = '  -1.91%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.54'
$ws.Range("E32").Value = '  -4.14%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0631'
$ws.Range("E33").Value = '  -5.41%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.36'
$ws.Range("E34").Value = '  -1.89%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.33'
$ws.Range("E35").Value = '  -8.69%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.50'
$ws.Range("E36").Value = '  -3.00%  '
$ws.Range("E37").Value = '  +1.05%  '
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.47'
$ws.Range("E39").Value = '  +2.46%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.02'
$ws.Range("E40").Value = '  +0.90%  '
$ws.Range("E41").Value = '  -0.50%  '
$ws.Range("D42").Value = '1.439.29'
$ws.Range("E42").Value = '  +5.13%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0208'
$ws.Range("E43").Value = '  -3.58%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0909'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '88.43'
$ws.Range("E45").Value = '  -2.72%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '15.42'
$ws.Range("E46").Value = '  -6.48%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.01'
$ws.Range("E47").Value = '  -3.01%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.89'
$ws.Range("E48").Value = '  +0.64%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.83'
$ws.Range("E49").Value = '  -7.90%  '
$ws.Range("D50").Value = '2.185.04'
$ws.Range("E50").Value = '  -0.68%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.92'
$ws.Range("E51").Value = '  -10.09%  '
